# Update the PSSM score matrix (B2:K21) with the supplemental-figure values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per residue (F, W, Y, P, M, I, L, V, A, G, C, S, T, N, Q, D, E, H, K, R),
# columns B:K hold the score for offsets 0-9.
$pssmRows = @(
    @(-18.63698210657446, 2.407260576718049, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # F
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # W
    @(-18.63698210657446, 2.130769956316451, 2.837280846297119, -18.63698210657446, 2.463857935638865, -18.63698210657446, 1.898317004384056, -18.63698210657446, 2.224060108389656, -18.63698210657446),  # Y
    @(-18.63698210657446, 1.063818253793698, -18.63698210657446, -18.63698210657446, -18.63698210657446, 2.156875969550921, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # P
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # M
    @(2.972767016741356, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # I
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, 2.857286080148574, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # L
    @(3.602831118432168, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # V
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, 4.321924732807638, -18.63698210657446, 2.187654279322414),  # A
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, 2.027542434368855, -18.63698210657446, 2.552044307456337, -18.63698210657446, -18.63698210657446, -18.63698210657446, 1.303849063223265),  # G
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # C
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, 1.680140089298963, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, 2.458935747896004, 1.558054814072276),  # S
    @(-18.63698210657446, -18.63698210657446, 1.691207408875055, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, 2.165279438257828),  # T
    @(-18.63698210657446, -18.63698210657446, -0.2331503757964419, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # N
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, 2.467548184711234, -18.63698210657446),  # Q
    @(-18.63698210657446, 0.7278899384779253, -0.01169777391712341, -18.63698210657446, -18.63698210657446, -18.63698210657446, 0.6968305721743945, -18.63698210657446, 1.093078326974468, -18.63698210657446),  # D
    @(-18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, -18.63698210657446, 0.6451312331230463, -18.63698210657446, 1.114681261243513, -18.63698210657446),  # E
    @(-18.63698210657446, -18.63698210657446, 1.708311594837235, -18.63698210657446, -18.63698210657446, -18.63698210657446, 1.91798608163757, -18.63698210657446, -18.63698210657446, -18.63698210657446),  # H
    @(-18.63698210657446, 1.640914809880734, 2.174245967960193, -18.63698210657446, 3.856294089964901, -18.63698210657446, 2.306640894559942, -18.63698210657446, -18.63698210657446, 2.471280075254948),  # K
    @(-18.63698210657446, 1.786799610829001, -18.63698210657446, 2.452006000987414, -18.63698210657446, 3.274385375586912, 2.124696536404199, -18.63698210657446, -18.63698210657446, -18.63698210657446)  # R
)

$numRows = $pssmRows.Count
$numCols = $pssmRows[0].Count
$values = New-Object "object[,]" $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $values[$r, $c] = $pssmRows[$r][$c]
    }
}

$ws.Range("B2:K21").Value = $values

